{"js": "// The document has a stray empty paragraph (indented/justified as a body\n// paragraph) sitting right before the paragraph that holds the \"Figura 1\"\n// architecture-diagram picture. That blank paragraph is leftover clutter\n// and gets removed, so the image paragraph follows immediately after the\n// \"Na Figura 1 pode ser visualizada...\" sentence.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Load inline pictures for every paragraph that follows an empty one so we\n// can recognise the blank paragraph that immediately precedes the figure.\nfor (let i = 0; i < items.length - 1; i++) {\n  if (items[i].text.trim().length === 0) {\n    items[i + 1].inlinePictures.load(\"items\");\n  }\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < items.length - 1; i++) {\n  if (items[i].text.trim().length === 0) {\n    const next = items[i + 1];\n    if (next.inlinePictures.items.length > 0) {\n      target = items[i];\n      break;\n    }\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# The document contains a stray empty paragraph (indented as if it were a\n# body paragraph, i.e. w:ind/@firstLine set and justified) sitting right\n# before the paragraph that holds \"Figura 1\" (the architecture diagram\n# image). That blank paragraph is pure clutter left over from editing and\n# is removed, pulling the image paragraph directly after the\n# \"Na Figura 1 pode ser visualizada...\" sentence.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$target = $null\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n\n    # Looking for an empty paragraph (no visible text besides the pilcrow)\n    if ($p.Range.Text.Trim().Length -eq 0) {\n        if ($i -lt $count) {\n            $nextP = $d.Paragraphs.Item($i + 1)\n            # ... immediately followed by the paragraph that holds the figure image\n            if ($nextP.Range.InlineShapes.Count -gt 0) {\n                $target = $p\n                break\n            }\n        }\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
